$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- New summary formulas -------------------------------------------------

# Row 12: average of column J (|S*|/n) under the data table
$ws.Range("J12").Formula = "=AVERAGE(J2:J11)"

# Row 14-17: labeled aggregate statistics
$ws.Range("A14").Value = "Average of SW(S*)/SW(OPT)"
$ws.Range("B14").Formula = "=AVERAGE(N2:N11)"

$ws.Range("A15").Value = "Average of SC(S*)/SC(OPT)"
$ws.Range("B15").Formula = "=AVERAGE(Z2:Z11)"

$ws.Range("A16").Value = "Worst of SW(S*)/SW(OPT)"
$ws.Range("B16").Formula = "=MIN(N2:N11)"

$ws.Range("A17").Value = "Worst of SC(S*)/SC(OPT)"
$ws.Range("B17").Formula = "=MAX(Z2:Z11)"

# ---- Formatting ------------------------------------------------------------

# J12 -> bold
$f1 = $ws.Range("J12").Font
$f1.Bold = $true
$f1.ThemeColor = 1

# B14:B17 -> bold, size 12, vertically centered
$r = $ws.Range("B14:B17")
$f2 = $r.Font
$f2.Bold = $true
$f2.ThemeColor = 1
$f2.Size = 12
$r.VerticalAlignment = -4108

$ws.Range("A14:B17").RowHeight = 15.6

$ws.Range("A14:B17").Select()
